$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5943607454143283
$ws.Range("C2").Value = -0.7999902782140788
$ws.Range("D2").Value = 1.704534436060835
$ws.Range("E2").Value = -0.4942653685160635
$ws.Range("F2").Value = 0.9369121485761228
$ws.Range("G2").Value = 0.272830226748868
$ws.Range("H2").Value = 0.7671803625714035
$ws.Range("I2").Value = 0.1142017009502766
$ws.Range("J2").Value = 0.539033769963857
$ws.Range("K2").Value = 0.5940197511622507
$ws.Range("B3").Value = -0.7775622985810702
$ws.Range("C3").Value = 1.707137975623284
$ws.Range("D3").Value = -0.5163330541700787
$ws.Range("E3").Value = 0.9309301682717679
$ws.Range("F3").Value = 0.2674869765664569
$ws.Range("G3").Value = 0.7568118292010373
$ws.Range("H3").Value = 0.1057243788266528
$ws.Range("I3").Value = 0.5312614616570462
$ws.Range("J3").Value = 0.5853957776292821
$ws.Range("K3").Value = 0.7117978212943072
$ws.Range("B4").Value = 1.728562547080504
$ws.Range("C4").Value = -0.3961552053708663
$ws.Range("D4").Value = 0.8417358238579847
$ws.Range("E4").Value = 0.2377599727178791
$ws.Range("F4").Value = 0.7654379492309473
$ws.Range("G4").Value = 0.07978420540121761
$ws.Range("H4").Value = 0.5084828986288187
$ws.Range("I4").Value = 0.5717777424806643
$ws.Range("J4").Value = 0.6932770241091315
$ws.Range("K4").Value = -0.2007560336349775
$ws.Range("B5").Value = -0.4384758376912558
$ws.Range("C5").Value = 0.8125313500022515
$ws.Range("D5").Value = 0.2347393729129579
$ws.Range("E5").Value = 0.7465766656194669
$ws.Range("F5").Value = 0.05945307433601332
$ws.Range("G5").Value = 0.4939527431945709
$ws.Range("H5").Value = 0.5551366969883963
$ws.Range("I5").Value = 0.6756931786803813
$ws.Range("J5").Value = -0.2172299171650203
$ws.Range("K5").Value = 0.553751414566908
$ws.Range("B6").Value = 1.152729090620161
$ws.Range("C6").Value = 0.310110434696895
$ws.Range("D6").Value = 0.5554833793064679
$ws.Range("E6").Value = 0.08414698692274653
$ws.Range("F6").Value = 0.5033321467508669
$ws.Range("G6").Value = 0.5003258983770053
$ws.Range("H6").Value = 0.6567990877917116
$ws.Range("I6").Value = -0.2320939179800661
$ws.Range("J6").Value = 0.5250397039375373
$ws.Range("K6").Value = 0.2638965897873631
$ws.Range("B7").Value = 0.7608053066871455
$ws.Range("C7").Value = 0.6025940815188262
$ws.Range("D7").Value = -0.1561362438680617
$ws.Range("E7").Value = 0.5392772278167297
$ws.Range("F7").Value = 0.4982751409374124
$ws.Range("G7").Value = 0.5802308483556813
$ws.Range("H7").Value = -0.2594933123553494
$ws.Range("I7").Value = 0.4987506537398119
$ws.Range("J7").Value = 0.2204992990740305
$ws.Range("B8").Value = 0.9149196684423646
$ws.Range("C8").Value = -0.02284750413253739
$ws.Range("D8").Value = 0.3601988060005381
$ws.Range("E8").Value = 0.5264693797079796
$ws.Range("F8").Value = 0.6162561595480749
$ws.Range("G8").Value = -0.2970348825595631
$ws.Range("H8").Value = 0.4932083146524507
$ws.Range("I8").Value = 0.2246746280127792
$ws.Range("B9").Value = 0.2127429869753038
$ws.Range("C9").Value = 0.4448775468748477
$ws.Range("D9").Value = 0.3808690130742625
$ws.Range("E9").Value = 0.6262099024073174
$ws.Range("F9").Value = -0.2816697128600181
$ws.Range("G9").Value = 0.4543599796950276
$ws.Range("H9").Value = 0.2088288189855932
$ws.Range("B10").Value = 0.7559017333562305
$ws.Range("C10").Value = 0.4979577874854577
$ws.Range("D10").Value = 0.4637457609577506
$ws.Range("E10").Value = -0.2527107623948165
$ws.Range("F10").Value = 0.4901028521499312
$ws.Range("G10").Value = 0.1775011726019661
$ws.Range("B11").Value = 0.7447829648895721
$ws.Range("C11").Value = 0.4814016284956401
$ws.Range("D11").Value = -0.347333001505811
$ws.Range("E11").Value = 0.5223202403984113
$ws.Range("F11").Value = 0.1898892984296834
$ws.Range("B12").Value = 0.7210779879118521
$ws.Range("C12").Value = -0.2623087580365975
$ws.Range("D12").Value = 0.4058356620403972
$ws.Range("E12").Value = 0.2049945700815359
$ws.Range("B13").Value = -0.09744868100251025
$ws.Range("C13").Value = 0.4193729342883134
$ws.Range("D13").Value = 0.1420216510915729
$ws.Range("B14").Value = 0.6732219761537215
$ws.Range("C14").Value = 0.2413397012736094
$ws.Range("B15").Value = 0.2853993925130583

$ws.Range("K7").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("B16").ClearContents()
